$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.541.63"
$ws.Range("E2").Value = "  +4.12%  "

$ws.Range("D3").Value = "4.030.04"
$ws.Range("E3").Value = "  +4.13%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'531.74"
$ws.Range("E5").Value = "  +1.68%  "

$ws.Range("D6").Value = "'150.02"
$ws.Range("E6").Value = "  +6.24%  "

$ws.Range("D7").Value = "'0.689"
$ws.Range("E7").Value = "  +13.32%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.759"
$ws.Range("E9").Value = "  +6.42%  "

$ws.Range("D10").Value = "'0.174"
$ws.Range("E10").Value = "  +4.29%  "

$ws.Range("D11").Value = "'0.0000330"
$ws.Range("E11").Value = "  +3.38%  "

$ws.Range("D12").Value = "'48.04"
$ws.Range("E12").Value = "  +15.38%  "

$ws.Range("D13").Value = "'10.85"
$ws.Range("E13").Value = "  +4.76%  "

$ws.Range("D14").Value = "4.671.50"
$ws.Range("E14").Value = "  +4.14%  "

$ws.Range("D15").Value = "4.017.97"
$ws.Range("E15").Value = "  +3.90%  "

$ws.Range("D16").Value = "'14.29"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").Value = "'20.70"
$ws.Range("E17").Value = "  -2.74%  "

$ws.Range("D18").Value = "'1.21"
$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "71.458.65"
$ws.Range("E20").Value = "  +4.02%  "

$ws.Range("D21").Value = "'433.18"
$ws.Range("E21").Value = "  +3.89%  "

$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "'99.09"
$ws.Range("E22").Value = "  +14.09%  "

$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'3.63"
$ws.Range("E23").Value = "  +3.09%  "

$ws.Range("D24").Value = "'14.65"
$ws.Range("E24").Value = "  +4.68%  "

$ws.Range("D25").Value = "'4.23"
$ws.Range("E25").Value = "  +6.51%  "

$ws.Range("D26").Value = "'11.27"
$ws.Range("E26").Value = "  -1.34%  "

$ws.Range("D27").Value = "'10.79"
$ws.Range("E27").Value = "  +2.57%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'5.83"
$ws.Range("E28").Value = "  +2.92%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'36.98"
$ws.Range("E29").Value = "  +4.10%  "

$ws.Range("D30").Value = "'3.49"
$ws.Range("E30").Value = "  +23.40%  "

$ws.Range("D31").Value = "'13.53"
$ws.Range("E31").Value = "  +2.24%  "

$ws.Range("D32").Value = "'0.131"
$ws.Range("E32").Value = "  +4.96%  "

$ws.Range("D33").Value = "'678.80"
$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("D34").Value = "'6.86"
$ws.Range("E34").Value = "  +0.36%  "

$ws.Range("D35").Value = "'66.50"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").Value = "'42.40"
$ws.Range("E36").Value = "  +6.88%  "

$ws.Range("E37").Value = "  -3.64%  "

$ws.Range("D38").Value = "'0.157"
$ws.Range("E38").Value = "  +5.76%  "

$ws.Range("D39").Value = "0.0₃0841"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("D40").Value = "'3.49"
$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").Value = "'0.0493"
$ws.Range("E42").Value = "  +4.00%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "'3.19"
$ws.Range("E44").Value = "  +2.12%  "

$ws.Range("E45").Value = "  +8.38%  "

$ws.Range("D46").Value = "'2.71"
$ws.Range("E46").Value = "  -4.16%  "

$ws.Range("D47").Value = "'3.43"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("D48").Value = "'9.51"
$ws.Range("E48").Value = "  +10.45%  "

$ws.Range("D49").Value = "'3.03"
$ws.Range("E49").Value = "  +0.65%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000274"
$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.35"
$ws.Range("E51").Value = "  +2.00%  "
